# Updated cryptos list on Tue Jul 16 22:16:49 UTC 2024 with GitHub Actions
# Refresh price/volume(1h) values and swap dogwifhat/Cosmos rows (50 <-> 51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "64.328.01"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value2 = "  +0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.422.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value2 = "  -0.71%  "
$ws.Range("E4").Value2 = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "572.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "159.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +1.61%  "
$ws.Range("E7").Value2 = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.420.89"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.579"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  +8.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.31"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  -3.74%  "
$ws.Range("E11").Value2 = "  +0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.436"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  -0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.010.81"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  -0.37%  "
$ws.Range("E14").Value2 = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0000192"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  +2.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "28.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  +1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "64.380.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.434.02"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.32"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "14.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "383.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -2.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "8.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  -4.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "72.92"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.541"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -0.35%  "
$ws.Range("E25").Value2 = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.0000121"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  +11.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "9.59"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.179"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -1.15%  "
$ws.Range("E29").Value2 = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "6.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  +5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  +1.45%  "
$ws.Range("E32").Value2 = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "6.52"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "23.47"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "7.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  +2.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "162.50"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  +2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "1.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "3.021.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  +5.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0761"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = "  -2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "27.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -3.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "4.51"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  +2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "42.74"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  +1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.0316"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.766"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "24.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  +8.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.875"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  +5.03%  "

# Rows 50 and 51 swapped position with refreshed data:
# previously row50 = dogwifhat, row51 = Cosmos -> now row50 = Cosmos, row51 = dogwifhat
$ws.Range("B50").Value2 = "Cosmos"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "6.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  +2.96%  "

$ws.Range("B51").Value2 = "dogwifhat"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = "  +3.12%  "
